# Applies the "fedcore" -> "approach" anonymization, the matching border
# restyle on the merged-header spacer cells, and drops the stray empty
# G5 cell on the computational_comparison sheet.

$wb = $excel.ActiveWorkbook

function Set-SpacerBorder($ws, $cellAddr, $removeLeft, $removeRight) {
    $rng = $ws.Range($cellAddr)
    # Reset to the default (un-styled) cell style first, so the new style
    # is built cleanly from xf 0 (plain font, no alignment) instead of
    # inheriting the bold/centered header style.
    $rng.Style = "Normal"
    # Re-apply a thin top+bottom box (matches the original boxed style),
    # then strip the sides that should no longer show a border.
    $rng.Borders(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders(9).LineStyle = 1   # xlEdgeBottom
    if ($removeLeft) {
        $rng.Borders(7).LineStyle = -4142   # xlEdgeLeft -> xlLineStyleNone
    } else {
        $rng.Borders(7).LineStyle = 1
    }
    if ($removeRight) {
        $rng.Borders(10).LineStyle = -4142  # xlEdgeRight -> xlLineStyleNone
    } else {
        $rng.Borders(10).LineStyle = 1
    }
}

# ---- Sheet 1: quality_comparison ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-SpacerBorder $ws1 "C1" $true  $true
Set-SpacerBorder $ws1 "D1" $true  $false

$ws1.Range("C2").Value = "approach"

# ---- Sheet 2: computational_comparison ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-SpacerBorder $ws2 "C1" $true  $true
Set-SpacerBorder $ws2 "D1" $true  $false
Set-SpacerBorder $ws2 "F1" $true  $true
Set-SpacerBorder $ws2 "G1" $true  $false

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell entirely.
$ws2.Range("G5").ClearContents()
